$d = $word.ActiveDocument

# Locate the Heading1 paragraph "Cadena de suministros" that is immediately
# followed by the Heading2 paragraph "Conceptos básicos" (the third
# occurrence of the heading text in this document).
$conceptosPara = $null
foreach ($p in $d.Paragraphs) {
    $txt = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq "Cadena de suministros") {
        $next = $p.Next()
        if ($next -ne $null) {
            $nextTxt = $next.Range.Text.TrimEnd([char]13, [char]7)
            if ($nextTxt -eq "Conceptos básicos") {
                $conceptosPara = $next
            }
        }
    }
}

if ($conceptosPara -eq $null) {
    throw "Could not find 'Conceptos básicos' paragraph following 'Cadena de suministros'"
}

# Remember the preceding paragraph so we can fetch the freshly inserted one.
$precedingPara = $conceptosPara.Previous()

# Insert a new empty paragraph right before "Conceptos básicos".
$conceptosPara.Range.InsertParagraphBefore()

$newPara = $precedingPara.Next()
$newPara.Range.Text = "hola como estas"
$newPara.Style = "FirstParagraph"
